$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.644.49"
$ws.Range("E2").Value = "  +5.91%  "
# Row 3
$ws.Range("D3").Value = "2.065.04"
$ws.Range("E3").Value = "  +4.38%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").Value = "'253.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.52%  "
# Row 6
$ws.Range("D6").Value = "'0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.53%  "
# Row 7
$ws.Range("D7").Value = "'65.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.05%  "
# Row 8
$ws.Range("E8").Value = "  -0.03%  "
# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.73%  "
# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'59.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
# Row 11
$ws.Range("D11").Value = "'0.0771"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.45%  "
# Row 12
$ws.Range("E12").Value = "  +1.43%  "
# Row 13
$ws.Range("D13").Value = "'0.910"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.06%  "
# Row 14
$ws.Range("D14").Value = "'14.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.51%  "
# Row 15
$ws.Range("D15").Value = "2.368.02"
$ws.Range("E15").Value = "  +4.31%  "
# Row 16
$ws.Range("D16").Value = "'5.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.95%  "
# Row 17
$ws.Range("D17").Value = "'20.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +18.59%  "
# Row 18
$ws.Range("D18").Value = "2.065.64"
$ws.Range("E18").Value = "  +4.45%  "
# Row 19
$ws.Range("D19").Value = "37.544.89"
$ws.Range("E19").Value = "  +5.87%  "
# Row 20
$ws.Range("D20").Value = "'74.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.47%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0882"
$ws.Range("E21").Value = "  +4.45%  "
# Row 22
$ws.Range("D22").Value = "'5.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.36%  "
# Row 23
$ws.Range("D23").Value = "'241.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.94%  "
# Row 24
$ws.Range("D24").Value = "'2.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.03%  "
# Row 25
$ws.Range("E25").Value = "  +0.03%  "
# Row 26
$ws.Range("D26").Value = "'2.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.51%  "
# Row 27
$ws.Range("D27").Value = "'9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.97%  "
# Row 28
$ws.Range("D28").Value = "'162.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "
# Row 29
$ws.Range("D29").Value = "'20.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.37%  "
# Row 30
$ws.Range("D30").Value = "'5.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.18%  "
# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.00%  "
# Row 32
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.113"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.45%  "
# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.15%  "
# Row 34
$ws.Range("E34").Value = "  +11.34%  "
# Row 35
$ws.Range("D35").Value = "'0.0623"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.49%  "
# Row 36
$ws.Range("D36").Value = "'2.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.74%  "
# Row 37
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +22.96%  "
# Row 38
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
# Row 39
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.81%  "
# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +33.81%  "
# Row 41
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.56%  "
# Row 42
$ws.Range("E42").Value = "  +3.98%  "
# Row 43
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.87%  "
# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.30%  "
# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0220"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.60%  "
# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'17.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.83%  "
# Row 47
$ws.Range("D47").Value = "'95.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.92%  "
# Row 48
$ws.Range("D48").Value = "'7.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.59%  "
# Row 49
$ws.Range("D49").Value = "1.423.24"
$ws.Range("E49").Value = "  +2.70%  "
# Row 50
$ws.Range("D50").Value = "'2.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.53%  "
# Row 51
$ws.Range("D51").Value = "'46.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.63%  "
